# Update the "Fitness" values in column C to reflect the new run results.
# The values change in blocks according to the "Generation" (column B) ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C66").Value = 8120
$ws.Range("C67:C74").Value = 8088
$ws.Range("C75:C88").Value = 8020
$ws.Range("C89:C97").Value = 7703
$ws.Range("C98:C252").Value = 7573
